# ---------------------------------------------------------------------------
# Update "update code tao report luong tai report co so":
#   1. Repurpose the existing "Lương" sheet as a new "Đơn 1 bác sĩ" sheet
#      (keeps sheetId=2 / rId2) and fill it with the per-service detail
#      rows used to compute doctor-1 commission ("Chiết khấu bác sĩ 1").
#   2. Insert a brand-new "Lương" sheet right after it (sheetId=3 / rId3)
#      and recompute the payroll summary rows with the figures derived
#      from the new "Đơn 1 bác sĩ" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Step 1: repurpose the old "Lương" worksheet -> "Đơn 1 bác sĩ" --------
$donBacSi = $wb.Worksheets.Item("Lương")
$donBacSi.Cells.Clear()
$donBacSi.Name = "Đơn 1 bác sĩ"

$headers = @(
    "Tiền tố",
    "Mã dịch vụ",
    "Ngày thực hiện",
    "Cơ sở",
    "Khách hàng",
    "Nguồn khách",
    "Tên dịch vụ",
    "Đơn giá gốc",
    "Sale phụ",
    "Upsale",
    "Đơn giá",
    "Đã thanh toán",
    "Tỉ lệ chiết khấu bác sĩ 1",
    "Chiết khấu bác sĩ 1"
)
for ($col = 1; $col -le $headers.Length; $col++) {
    $donBacSi.Cells.Item(1, $col).Value = $headers[$col - 1]
}

$donBacSi.Cells.Item(2, 1).Value = "HD-LUXURY"
$donBacSi.Cells.Item(2, 2).Value = 619
$donBacSi.Cells.Item(2, 3).NumberFormat = "@"
$donBacSi.Cells.Item(2, 3).Value = "08-02-2024"
$donBacSi.Cells.Item(2, 4).Value = "SÓC TRĂNG"
$donBacSi.Cells.Item(2, 5).Value = "mai hồng nương"
$donBacSi.Cells.Item(2, 6).Value = "Cá nhân"
$donBacSi.Cells.Item(2, 7).Value = "Thu cánh mũi"
$donBacSi.Cells.Item(2, 8).Value = 8000000
$donBacSi.Cells.Item(2, 11).Value = 8000000
$donBacSi.Cells.Item(2, 12).Value = 6000000
$donBacSi.Cells.Item(2, 13).Value = 0.1
$donBacSi.Cells.Item(2, 14).Value = 600000

$donBacSi.Cells.Item(3, 1).Value = "Tổng"
$donBacSi.Cells.Item(3, 2).Value = 1
$donBacSi.Cells.Item(3, 8).Value = 8000000
$donBacSi.Cells.Item(3, 10).Value = 0
$donBacSi.Cells.Item(3, 11).Value = 8000000
$donBacSi.Cells.Item(3, 12).Value = 6000000
$donBacSi.Cells.Item(3, 13).Value = 0
$donBacSi.Cells.Item(3, 14).Value = 600000

# --- Step 2: add a fresh "Lương" worksheet right after it -----------------
$luong = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $donBacSi)
$luong.Name = "Lương"

$payroll = @(
    @(1, 'Danh mục lương', 6),
    @(2, 'Tổng công tại CẦN THƠ', 2),
    @(3, 'Phụ cấp tại CẦN THƠ', 70000),
    @(4, 'Lương cơ bản tại CẦN THƠ', 238095.2380952381),
    @(5, 'Chiết khấu sale chính tại CẦN THƠ', 0),
    @(6, 'Chiết khấu sale phụ tại CẦN THƠ', 0),
    @(7, 'Đơn 1 bác sĩ tại CẦN THƠ', 0),
    @(8, 'Đơn 2 bác sĩ tại CẦN THƠ', 0),
    @(9, 'Công phụ phẫu 1 tại CẦN THƠ', 0),
    @(10, 'Công phụ phẫu 2 tại CẦN THƠ', 0),
    @(11, 'Ứng lương tại CẦN THƠ', 0),
    @(12, 'Tổng công tại LONG XUYÊN', 0),
    @(13, 'Lương công tác tại LONG XUYÊN', 0),
    @(14, 'Lương cơ bản tại LONG XUYÊN', 238095.2380952381),
    @(15, 'Chiết khấu sale chính tại LONG XUYÊN', 0),
    @(16, 'Chiết khấu sale phụ tại LONG XUYÊN', 0),
    @(17, 'Đơn 1 bác sĩ tại LONG XUYÊN', 0),
    @(18, 'Đơn 2 bác sĩ tại LONG XUYÊN', 0),
    @(19, 'Công phụ phẫu 1 tại LONG XUYÊN', 0),
    @(20, 'Công phụ phẫu 2 tại LONG XUYÊN', 0),
    @(21, 'Ứng lương tại LONG XUYÊN', 0),
    @(22, 'Tổng công tại SÓC TRĂNG', 0),
    @(23, 'Lương công tác tại SÓC TRĂNG', 0),
    @(24, 'Lương cơ bản tại SÓC TRĂNG', 238095.2380952381),
    @(25, 'Chiết khấu sale chính tại SÓC TRĂNG', 0),
    @(26, 'Chiết khấu sale phụ tại SÓC TRĂNG', 0),
    @(27, 'Đơn 1 bác sĩ tại SÓC TRĂNG', 600000),
    @(28, 'Đơn 2 bác sĩ tại SÓC TRĂNG', 0),
    @(29, 'Công phụ phẫu 1 tại SÓC TRĂNG', 0),
    @(30, 'Công phụ phẫu 2 tại SÓC TRĂNG', 0),
    @(31, 'Ứng lương tại SÓC TRĂNG', 0),
    @(32, 'Tổng lương tại CẦN THƠ', 308095.2380952381),
    @(33, 'Tổng lương tại LONG XUYÊN', 238095.2380952381),
    @(34, 'Tổng lương tại SÓC TRĂNG', 838095.2380952381),
    @(35, 'Tổng lương tại HỆ THỐNG', 1384285.714285714)
)

foreach ($item in $payroll) {
    $r = $item[0]
    $luong.Cells.Item($r, 1).Value = $item[1]
    $luong.Cells.Item($r, 2).Value = $item[2]
}

$wb.Worksheets.Item("Đơn sale chính").Activate()
